$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 2.742348
$ws.Cells.Item(2, 8).Value = 8.227044
$ws.Cells.Item(2, 9).Value = 0.8112951562181048
$ws.Cells.Item(2, 10).Value = 0.8112951562181048
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 20.98493066666667
$ws.Cells.Item(2, 14).Value = 62.954792
$ws.Cells.Item(2, 15).Value = 0.6651705782252746
$ws.Cells.Item(2, 16).Value = 0.6651705782252747
$ws.Cells.Item(2, 17).Value = 57.547982643872
$ws.Cells.Item(2, 18).Value = 517.931843794848
$ws.Cells.Item(2, 19).Value = 0.5396496681729612
$ws.Cells.Item(2, 20).Value = 0.5396496681729613

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 2.742348
$ws.Cells.Item(3, 8).Value = 8.227044
$ws.Cells.Item(3, 9).Value = 0.8112951562181048
$ws.Cells.Item(3, 10).Value = 0.8112951562181048
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 7.288088333333334
$ws.Cells.Item(3, 14).Value = 21.864265
$ws.Cells.Item(3, 15).Value = 0.231014436399387
$ws.Cells.Item(3, 16).Value = 0.231014436399387
$ws.Cells.Item(3, 17).Value = 19.98647446474
$ws.Cells.Item(3, 18).Value = 179.87827018266
$ws.Cells.Item(3, 19).Value = 0.1874208932672781
$ws.Cells.Item(3, 20).Value = 0.1874208932672781

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 2.742348
$ws.Cells.Item(4, 8).Value = 8.227044
$ws.Cells.Item(4, 9).Value = 0.8112951562181048
$ws.Cells.Item(4, 10).Value = 0.8112951562181048
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 3.275175333333333
$ws.Cells.Item(4, 14).Value = 9.825526
$ws.Cells.Item(4, 15).Value = 0.1038149853753384
$ws.Cells.Item(4, 16).Value = 0.1038149853753384
$ws.Cells.Item(4, 17).Value = 8.981670525016
$ws.Cells.Item(4, 18).Value = 80.835034725144
$ws.Cells.Item(4, 19).Value = 0.08422459477786541
$ws.Cells.Item(4, 20).Value = 0.08422459477786541

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.3936376666666666
$ws.Cells.Item(5, 8).Value = 1.180913
$ws.Cells.Item(5, 9).Value = 0.1164536128425946
$ws.Cells.Item(5, 10).Value = 0.1164536128425946
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 20.98493066666667
$ws.Cells.Item(5, 14).Value = 62.954792
$ws.Cells.Item(5, 15).Value = 0.6651705782252746
$ws.Cells.Item(5, 16).Value = 0.6651705782252747
$ws.Cells.Item(5, 17).Value = 8.260459142788443
$ws.Cells.Item(5, 18).Value = 74.344132285096
$ws.Cells.Item(5, 19).Value = 0.0774615169909309
$ws.Cells.Item(5, 20).Value = 0.07746151699093091

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.3936376666666666
$ws.Cells.Item(6, 8).Value = 1.180913
$ws.Cells.Item(6, 9).Value = 0.1164536128425946
$ws.Cells.Item(6, 10).Value = 0.1164536128425946
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 7.288088333333334
$ws.Cells.Item(6, 14).Value = 21.864265
$ws.Cells.Item(6, 15).Value = 0.231014436399387
$ws.Cells.Item(6, 16).Value = 0.231014436399387
$ws.Cells.Item(6, 17).Value = 2.868866085993889
$ws.Cells.Item(6, 18).Value = 25.819794773945
$ws.Cells.Item(6, 19).Value = 0.0269024657375044
$ws.Cells.Item(6, 20).Value = 0.0269024657375044

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.3936376666666666
$ws.Cells.Item(7, 8).Value = 1.180913
$ws.Cells.Item(7, 9).Value = 0.1164536128425946
$ws.Cells.Item(7, 10).Value = 0.1164536128425946
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 3.275175333333333
$ws.Cells.Item(7, 14).Value = 9.825526
$ws.Cells.Item(7, 15).Value = 0.1038149853753384
$ws.Cells.Item(7, 16).Value = 0.1038149853753384
$ws.Cells.Item(7, 17).Value = 1.289232376137555
$ws.Cells.Item(7, 18).Value = 11.603091385238
$ws.Cells.Item(7, 19).Value = 0.01208963011415927
$ws.Cells.Item(7, 20).Value = 0.01208963011415928

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.2442243333333333
$ws.Cells.Item(8, 8).Value = 0.732673
$ws.Cells.Item(8, 9).Value = 0.07225123093930062
$ws.Cells.Item(8, 10).Value = 0.07225123093930062
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 20.98493066666667
$ws.Cells.Item(8, 14).Value = 62.954792
$ws.Cells.Item(8, 15).Value = 0.6651705782252746
$ws.Cells.Item(8, 16).Value = 0.6651705782252747
$ws.Cells.Item(8, 17).Value = 5.125030702112889
$ws.Cells.Item(8, 18).Value = 46.125276319016
$ws.Cells.Item(8, 19).Value = 0.04805939306138245
$ws.Cells.Item(8, 20).Value = 0.04805939306138245

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.2442243333333333
$ws.Cells.Item(9, 8).Value = 0.732673
$ws.Cells.Item(9, 9).Value = 0.07225123093930062
$ws.Cells.Item(9, 10).Value = 0.07225123093930062
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 7.288088333333334
$ws.Cells.Item(9, 14).Value = 21.864265
$ws.Cells.Item(9, 15).Value = 0.231014436399387
$ws.Cells.Item(9, 16).Value = 0.231014436399387
$ws.Cells.Item(9, 17).Value = 1.779928514482778
$ws.Cells.Item(9, 18).Value = 16.019356630345
$ws.Cells.Item(9, 19).Value = 0.01669107739460448
$ws.Cells.Item(9, 20).Value = 0.01669107739460448

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.2442243333333333
$ws.Cells.Item(10, 8).Value = 0.732673
$ws.Cells.Item(10, 9).Value = 0.07225123093930062
$ws.Cells.Item(10, 10).Value = 0.07225123093930062
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 3.275175333333333
$ws.Cells.Item(10, 14).Value = 9.825526
$ws.Cells.Item(10, 15).Value = 0.1038149853753384
$ws.Cells.Item(10, 16).Value = 0.1038149853753384
$ws.Cells.Item(10, 17).Value = 0.7998775123331111
$ws.Cells.Item(10, 18).Value = 7.198897610998
$ws.Cells.Item(10, 19).Value = 0.007500760483313691
$ws.Cells.Item(10, 20).Value = 0.007500760483313691

Write-Output "done"